$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new sheet "SR03-JP" as the last sheet (after SPDS-JP)
# ---------------------------------------------------------------------------
$sheetSPDS = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $sheetSPDS)
$ws.Name = "SR03-JP"

# ---------------------------------------------------------------------------
# 2. Populate the new sheet's content.
#    Shared-string insertion order matters (it dictates the sst index), so
#    write the brand new strings in the same order they appear in the target
#    workbook: title, then the three card names (A column top-to-bottom),
#    then the trailing ";" separator string (E column).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Structure Deck R: Machine Dragon Re-Volt"

$ws.Range("A2").Value = "Ancient Gear Gadget"
$ws.Range("A3").Value = "Ancient Gear Reactor Dragon"
$ws.Range("A4").Value = "Ancient Gear Catapult"

$ws.Range("B2").Value = 100303000
$ws.Range("B3").Value = 100303001
$ws.Range("B4").Value = 100303021

$ws.Range("C2").Value = ":"
$ws.Range("C3").Value = ":"
$ws.Range("C4").Value = ":"

$ws.Range("E2").Value = ";"
$ws.Range("E3").Value = ";"
$ws.Range("E4").Value = ";"

# Date value (2016-09-24, Excel serial 42637) for the header row.
$ws.Range("B1").Value = 42637

# ---------------------------------------------------------------------------
# 3. Formatting: copy the header styles (bold/bordered title cell + the
#    bordered date cell) from an existing sheet so we reuse the workbook's
#    existing style/border entries instead of minting new ones, then strip
#    the bold weight off the date cell to match the new header style.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(3)

$templateSheet.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$templateSheet.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B1").Font.Bold = $false

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Column widths / row height to roughly match the source sheets.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 37.07421875
$ws.Columns.Item(2).ColumnWidth = 9.84375
$ws.Columns.Item(3).ColumnWidth = 1.3828125
$ws.Columns.Item(5).ColumnWidth = 1.3828125
$ws.Rows.Item(1).RowHeight = 24

# Match the portrait page setup used by the other sheets.
$ws.PageSetup.Orientation = 1   # xlPortrait

# ---------------------------------------------------------------------------
# 5. Selections: update the active cell on every sheet per the target, and
#    make sure the previously-active sheet (SPDS-JP) stays the active tab.
# ---------------------------------------------------------------------------
$ws.Range("B8").Select() | Out-Null

$sheetVP16 = $wb.Worksheets.Item(1)
$sheetVP16.Range("F8").Select() | Out-Null

$sheetTDIL = $wb.Worksheets.Item(2)
$sheetTDIL.Range("A16").Select() | Out-Null

$sheetSPDS.Activate()
$sheetSPDS.Range("B34").Select() | Out-Null
